# icStars development guide - "high level coding standards" update
#
# Four changes, per the commit's diff:
#  1. Remove the "_GoBack" bookmark that used to sit right after the
#     apple.png image path text.
#  2. Insert a new CSS bullet ("All styles should reside in CSS files
#     located in the /css folder") immediately before the existing
#     "Ensure that styles are properly indented" bullet, using the same
#     list (ListParagraph / numId 2) formatting.
#  3. Re-add the "_GoBack" bookmark (collapsed, i.e. zero-length) right
#     after the text of the "Be consistent" bullet.
#  4. Drop the trailing period from " if you need to." so it reads
#     " if you need to".

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark (after "apple.png") ----------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Insert the new "All styles..." bullet ------------------------------
#        immediately before the "Ensure that styles are properly indented"
#        bullet.
$find = $d.Content.Find
$found = $find.Execute("Ensure that styles are properly indented", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetPara = $find.Parent.Duplicate.Paragraphs(1)
    $newParaRange = $targetPara.Range

    # InsertParagraphBefore copies the pPr (ListParagraph style, numId 2
    # bullet) of the paragraph it's called on into a new, still-empty
    # paragraph inserted immediately before it, and repoints $newParaRange
    # at that freshly created paragraph.
    $newParaRange.InsertParagraphBefore()
    $newParaRange.Text = "All styles should reside in CSS files located in the /css folder"
}

# --- 3. Re-add the _GoBack bookmark right after "Be consistent" ------------
$find2 = $d.Content.Find
$found2 = $find2.Execute("Be consistent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $afterConsistent = $find2.Parent.Duplicate
    $afterConsistent.Collapse(0)

    # Bookmarks.Add doesn't reliably anchor a truly zero-length Range in
    # this object model, so insert a one-character placeholder, bookmark
    # across that character, then clear the placeholder's text. The
    # bookmark itself stays put, now collapsed in the right spot.
    $afterConsistent.InsertAfter("#")
    $bmRange = $d.Range($afterConsistent.Start, $afterConsistent.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    $d.Bookmarks("_GoBack").Range.Text = ""
}

# --- 4. " if you need to." -> " if you need to" ----------------------------
[void]$d.Content.Find.Execute(" if you need to.", $true, $false, $false, $false, $false, $true, 1, $false, " if you need to", 2)
